$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-27 all carry the same serial date value
# (45192 -> 2023-09-23). The update bumps them all to 45202 (2023-10-03).
foreach ($row in 2..27) {
    $ws.Cells.Item($row, 3).Value = 45202
}
